# Daily update at 8 AM UTC
# Appends the next day's row to the "Wins Over Time" tracking sheet and
# moves the "last row" date formatting from the previous last row to the
# newly appended one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new day's data as the new last row (row 18).
$ws.Range("A18").Value = 45967
$ws.Range("B18").Value = 36
$ws.Range("C18").Value = 42
$ws.Range("D18").Value = 43

# The newly appended row becomes the sheet's last row, so it takes on the
# distinct "last row" date format that row 17 used to have.
$ws.Range("A18").NumberFormat = $ws.Range("A17").NumberFormat

# Row 17 is no longer the last row, so it reverts to the regular date
# format used by the rest of the data rows.
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
